# Insert a new data row at row 145 (pushing the existing rows 145-169 down
# to 146-170) on the active sheet, then populate the new row with the
# "Macroferia Regional de Talca" / Mandarina / Murcott record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 145, shifting rows 145:169 down
# to 146:170 (same as right-clicking the row header -> Insert in Excel).
$ws.Rows(145).Insert()

# Populate the newly inserted row 145 with the new record.
$ws.Range("A145").Value = 5
$ws.Range("B145").Value = "Macroferia Regional de Talca"
$ws.Range("C145").Value = "Maule"
$ws.Range("D145").Value = 44449
$ws.Range("E145").Value = 7
$ws.Range("F145").Value = "Fruta"
$ws.Range("G145").Value = 100102
$ws.Range("H145").Value = "Cítricos"
$ws.Range("I145").Value = 100102004
$ws.Range("J145").Value = "Mandarina"
$ws.Range("K145").Value = "Murcott"
$ws.Range("L145").Value = "Segunda"
$ws.Range("M145").Value = 200
$ws.Range("N145").Value = 5000
$ws.Range("O145").Value = 5000
$ws.Range("P145").Value = 5000
$ws.Range("Q145").Value = "$/caja 18 kilos"
$ws.Range("R145").Value = "Provincia de Quillota"
$ws.Range("S145").Value = 278
$ws.Range("T145").Value = 18
